# Updated cryptos list on Mon Oct 23 10:30:27 UTC 2023 with GitHub Actions
#
# The "Price" column (D) holds values that look numeric (e.g. "0.999",
# "219.76") but are stored as plain text in the source data (note values
# like "30.511.40" or "1.913.86" which aren't valid numbers at all, and
# the mixed formatting throughout the column). A leading apostrophe is
# used on every D-column write below to force Excel to keep/treat the
# entry as text instead of auto-converting number-looking strings into
# numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    $range.Value = "'" + $text
}

# --- Row 2: Bitcoin ---
Set-TextValue $ws.Range("D2") "30.511.40"
$ws.Range("E2").Value = "  +2.02%  "

# --- Row 3: Ethereum ---
Set-TextValue $ws.Range("D3") "1.673.26"
$ws.Range("E3").Value = "  +2.47%  "

# --- Row 4: TetherUSD ---
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.22%  "

# --- Row 5: BNB ---
Set-TextValue $ws.Range("D5") "219.76"
$ws.Range("E5").Value = "  +2.49%  "

# --- Row 6: XRP ---
Set-TextValue $ws.Range("D6") "0.529"
$ws.Range("E6").Value = "  +2.06%  "

# --- Row 7: USDC ---
Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.12%  "

# --- Row 8: Solana ---
Set-TextValue $ws.Range("D8") "29.75"
$ws.Range("E8").Value = "  +4.65%  "

# --- Row 9: Cardano ---
$ws.Range("E9").Value = "  +2.76%  "

# --- Row 10: Dogecoin ---
Set-TextValue $ws.Range("D10") "0.0639"
$ws.Range("E10").Value = "  +5.20%  "

# --- Row 11: TRON ---
Set-TextValue $ws.Range("D11") "0.0905"
$ws.Range("E11").Value = "  -0.63%  "

# --- Row 12: WrappedliquidstakedEther2.0 ---
Set-TextValue $ws.Range("D12") "1.913.86"
$ws.Range("E12").Value = "  +2.51%  "

# --- Row 13: WrappedEther ---
Set-TextValue $ws.Range("D13") "1.683.22"
$ws.Range("E13").Value = "  +3.22%  "

# --- Row 14: Polygon ---
Set-TextValue $ws.Range("D14") "0.613"
$ws.Range("E14").Value = "  +9.03%  "

# --- Row 15: Chainlink ---
Set-TextValue $ws.Range("D15") "10.21"
$ws.Range("E15").Value = "  +10.22%  "

# --- Row 16: Polkadot ---
Set-TextValue $ws.Range("D16") "3.98"
$ws.Range("E16").Value = "  +3.49%  "

# --- Row 17: WrappedBTC ---
Set-TextValue $ws.Range("D17") "30.523.41"
$ws.Range("E17").Value = "  +1.91%  "

# --- Row 18: Litecoin ---
Set-TextValue $ws.Range("D18") "66.26"
$ws.Range("E18").Value = "  +3.62%  "

# --- Row 19: BitcoinCash ---
Set-TextValue $ws.Range("D19") "242.71"
$ws.Range("E19").Value = "  +0.26%  "

# --- Row 20: ShibaInu ---
Set-TextValue $ws.Range("D20") "0.0₃0721"
$ws.Range("E20").Value = "  +2.90%  "

# --- Row 21: Dai ---
Set-TextValue $ws.Range("D21") "0.999"
$ws.Range("E21").Value = "  -0.03%  "

# --- Row 22: Uniswap ---
Set-TextValue $ws.Range("D22") "4.26"
$ws.Range("E22").Value = "  +3.16%  "

# --- Row 23: Avalanche ---
Set-TextValue $ws.Range("D23") "9.99"
$ws.Range("E23").Value = "  +1.65%  "

# --- Row 24: Toncoin ---
Set-TextValue $ws.Range("D24") "2.15"
$ws.Range("E24").Value = "  +0.08%  "

# --- Row 25: Monero ---
Set-TextValue $ws.Range("D25") "158.24"
$ws.Range("E25").Value = "  +0.30%  "

# --- Row 26: EthereumClassic ---
Set-TextValue $ws.Range("D26") "15.86"
$ws.Range("E26").Value = "  +2.41%  "

# --- Row 27: Stellar ---
Set-TextValue $ws.Range("D27") "0.113"
$ws.Range("E27").Value = "  +2.45%  "

# --- Row 28: Cosmos ---
Set-TextValue $ws.Range("D28") "6.68"

# --- Row 29: BinanceUSD ---
Set-TextValue $ws.Range("D29") "0.999"
$ws.Range("E29").Value = "  -0.24%  "

# --- Row 30: Hedera ---
Set-TextValue $ws.Range("D30") "0.0495"
$ws.Range("E30").Value = "  +1.81%  "

# --- Row 31: PancakeSwap ---
Set-TextValue $ws.Range("D31") "1.15"
$ws.Range("E31").Value = "  +2.93%  "

# --- Row 32: Filecoin ---
Set-TextValue $ws.Range("D32") "3.47"
$ws.Range("E32").Value = "  +2.82%  "

# --- Row 33 & 34: InternetComputer(DFINITY)/Maker swapped places ---
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D33") "1.503.26"
$ws.Range("E33").Value = "  +5.39%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D34") "3.28"
$ws.Range("E34").Value = "  +3.41%  "

# --- Row 35: LidoDAOToken ---
$ws.Range("E35").Value = "  +6.82%  "

# --- Row 36: Aave ---
Set-TextValue $ws.Range("D36") "84.26"
$ws.Range("E36").Value = "  +10.99%  "

# --- Row 37: TrustWalletToken ---
$ws.Range("E37").Value = "  -0.84%  "

# --- Row 38: ImmutableX ---
Set-TextValue $ws.Range("D38") "0.598"
$ws.Range("E38").Value = "  +8.52%  "

# --- Row 39: VeChain ---
$ws.Range("E39").Value = "  +5.25%  "

# --- Row 40: MXToken ---
Set-TextValue $ws.Range("D40") "2.66"
$ws.Range("E40").Value = "  -4.88%  "

# --- Row 41: HuobiToken ---
$ws.Range("E41").Value = "  -0.38%  "

# --- Row 42: ARBITRUM ---
Set-TextValue $ws.Range("D42") "0.839"
$ws.Range("E42").Value = "  +1.66%  "

# --- Row 43: Kaspa ---
$ws.Range("E43").Value = "  +1.75%  "

# --- Row 44: RenderToken ---
$ws.Range("E44").Value = "  -1.06%  "

# --- Row 45: WEMIXToken ---
$ws.Range("E45").Value = "  +0.25%  "

# --- Row 46: PaxDollar ---
Set-TextValue $ws.Range("D46") "0.999"
$ws.Range("E46").Value = "  -0.12%  "

# --- Row 47 & 48: FraxShare/BitcoinSV swapped places ---
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws.Range("D47") "51.45"
$ws.Range("E47").Value = "  -2.83%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D48") "5.53"
$ws.Range("E48").Value = "  +3.39%  "

# --- Row 49: RocketPoolETH ---
Set-TextValue $ws.Range("D49") "1.804.10"
$ws.Range("E49").Value = "  +1.66%  "

# --- Row 50: Quant ---
Set-TextValue $ws.Range("D50") "94.47"
$ws.Range("E50").Value = "  +4.82%  "

# --- Row 51: BabyDogeCoin ---
$ws.Range("E51").Value = "  +0.71%  "
